$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new formatted row right after the current last data row (23),
# inheriting its cell formatting, then fill in the new skill entry:
# ID 22 "LimiterOverload" / StatusSkill / Mana 9 / CoolDown 3.
$ws.Rows("24").Insert(-4121) | Out-Null

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "LimiterOverload"
$ws.Range("C24").Value = "StatusSkill"
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 3

# Match the workbook's last-selection bookkeeping to the new last row
$ws.Range("A24:XFD24").Select() | Out-Null
